# PlayerPerformance_4609.xlsx edit script
# - Adds a new "Player Info" sheet (before "ODI Batting") with player bio data
# - Renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling",
#   replacing the full howstat URL values with just the numeric match code
# - Adds a new "ODI Batting Extra" sheet (after "ODI Bowling") with additional
#   per-match batting stats

$wb = $excel.ActiveWorkbook

$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$infoWs = $wb.Worksheets.Add($battingWs)
$infoWs.Name = "Player Info"

$infoHeader = $infoWs.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1

$infoWs.Range("A1").Value = "ID"
$infoWs.Range("B1").Value = "NAME"
$infoWs.Range("C1").Value = "BATTING_HAND"
$infoWs.Range("D1").Value = "BOWL_STYLE"

$infoWs.Range("A2").NumberFormat = "@"
$infoWs.Range("A2").Value = "4609"
$infoWs.Range("B2").Value = "Lachlan Hammond Ferguson"
$infoWs.Range("C2").Value = "Right Handed"
$infoWs.Range("D2").Value = "Right Arm Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" - rename MATCH_CARD_LINK column to MATCH_CODE and replace
#    the URL values with the bare numeric match code
# ---------------------------------------------------------------------------
# Re-fetch by name: inserting a sheet shifted worksheet positions, and this
# runtime's worksheet object references are position-based.
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

$battingLastRow = $battingWs.Cells.Item($battingWs.Rows.Count, 1).End(-4162).Row

$battingWs.Range("D1").Value = "MATCH_CODE"
$battingWs.Range("D2:D" + $battingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingWs.Cells.Item($r, 4)
    $url = $cell.Text
    $code = $url -replace '.*MatchCode=', ''
    $cell.Value = $code
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" - rename MATCH_CARD_LINK column to MATCH_CODE and replace
#    the URL values with the bare numeric match code
# ---------------------------------------------------------------------------
$bowlingLastRow = $bowlingWs.Cells.Item($bowlingWs.Rows.Count, 1).End(-4162).Row

$bowlingWs.Range("B1").Value = "MATCH_CODE"
$bowlingWs.Range("B2:B" + $bowlingLastRow).NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingWs.Cells.Item($r, 2)
    $url = $cell.Text
    $code = $url -replace '.*MatchCode=', ''
    $cell.Value = $code
}

# ---------------------------------------------------------------------------
# 4. "ODI Batting Extra" - new sheet inserted after "ODI Bowling"
# ---------------------------------------------------------------------------
# Re-fetch again in case anything shifted worksheet positions.
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$extraWs = $wb.Worksheets.Add($null, $bowlingWs)
$extraWs.Name = "ODI Batting Extra"

$extraHeader = $extraWs.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160
$extraHeader.Borders.LineStyle = 1

$extraWs.Range("A1").Value = "MATCH_CODE"
$extraWs.Range("B1").Value = "BATTING_POSITION"
$extraWs.Range("C1").Value = "NUM_4"
$extraWs.Range("D1").Value = "NUM_6"
$extraWs.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extraWs.Range("F1").Value = "MAN_OF_MATCH"

$extraWs.Range("A2:A21").NumberFormat = "@"
$extraWs.Range("C2:F21").NumberFormat = "@"

$extraData = @(
    @("4341", 10,   "0", "0", "",      "NO"),
    @("4353", $null, "", "", "",      "NO"),
    @("4355", $null, "", "", "",      "NO"),
    @("4423", 10,   "0", "0", "0.53%", "NO"),
    @("4605", 10,   "1", "0", "2.62%", "NO"),
    @("4614", 11,   "", "",  "",      "NO"),
    @("4625", 9,    "", "",  "",      "NO"),
    @("4636", 10,   "0", "0", "",      "NO"),
    @("4642", $null, "", "", "",      "NO"),
    @("4647", $null, "", "", "",      "NO"),
    @("4649", 10,   "0", "0", "1.65%", "NO"),
    @("4669", $null, "", "", "",      "NO"),
    @("4673", $null, "", "", "",      "NO"),
    @("4676", 11,   "0", "0", "",      "NO"),
    @("4686", $null, "", "", "",      "NO"),
    @("4688", $null, "", "", "",      "NO"),
    @("4690", 10,   "0", "0", "0.93%", "NO"),
    @("4692", 9,    "0", "0", "2.37%", "NO"),
    @("4695", $null, "", "", "",      ""),
    @("4697", $null, "", "", "",      "")
)

$r = 2
foreach ($row in $extraData) {
    $extraWs.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $extraWs.Cells.Item($r, 2).Value = $row[1]
    }
    $extraWs.Cells.Item($r, 3).Value = $row[2]
    $extraWs.Cells.Item($r, 4).Value = $row[3]
    $extraWs.Cells.Item($r, 5).Value = $row[4]
    $extraWs.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
